$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy the formatting from the row above (row 5, test case #3) down into row 6
# so the new row matches the look of the existing test-case rows.
$ws.Range("A5:E5").Copy()
$ws.Range("A6:E6").PasteSpecial(-4122)

# Fill in the new test case: "Navigate to Team Contributions Page"
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "test_<Navigate-To-TeamContributionsPage>"
$ws.Range("C6").Value = "This is to test whether users are able to navigate to Team Contributions page"
$ws.Range("D6").Value = "NIL"
$ws.Range("E6").Value = "Team contributions page is shown"

# Update the active selection to G7, matching the author's final cursor position
$ws.Range("G7").Select()
